$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Row, C, D, E, F, G, H, L, M, N (new loading_percent values)
$data = @(
    @(2, 11.37935798576697, 7.411053396066094, 11.3040552221086, 41.36510220724429, 59.94814049838629, 20.80522777188292, 8.220975979374128, 29.95001582830196, 17.0579469673359),
    @(3, 11.43552242938932, 7.427622296517418, 11.36945326088068, 40.77208391134781, 58.74362003227039, 20.69163462128778, 8.233131206538529, 29.11728923742223, 16.87438744843757),
    @(4, 11.47308690606094, 7.44024619031098, 11.41138837102974, 40.42137160076845, 58.01958400657601, 20.62907484195075, 8.241205883930579, 28.59824295286601, 16.76299461469746),
    @(5, 11.48916393077421, 7.446000449481871, 11.42892709687017, 40.28198847418093, 57.72887320278308, 20.6053974851434, 8.244650616712297, 28.38508346869886, 16.7179780484048),
    @(6, 11.49187981825253, 7.446992590660664, 11.43186661935711, 40.25906213378909, 57.68087544861616, 20.60157577648938, 8.24523194249573, 28.3495984653523, 16.71052720154735),
    @(7, 11.47330061866043, 7.440321333124284, 11.41162308026629, 40.41947730779496, 58.01564524344451, 20.62874815726554, 8.241251715619519, 28.5953744498604, 16.76238591809111),
    @(8, 11.39808127849035, 7.416254016165416, 11.3262362184176, 41.15794926421987, 59.52985605004546, 20.76457563476161, 8.22504053305787, 29.66467698008979, 16.99441112344212),
    @(9, 11.27524237836686, 7.388775602445202, 11.17282170791407, 42.7047843048516, 62.60342879943907, 21.08740957886812, 8.198077980852789, 31.68767647639156, 17.45777380137307),
    @(10, 11.20035619703844, 7.380970609362561, 11.06852073899712, 43.89073325593697, 64.89967585878045, 21.35808615779806, 8.18118006266927, 33.11372186092327, 17.80054858503573),
    @(11, 11.1696969832157, 7.380182280073917, 11.0228678970385, 44.4386643487431, 65.94735690621984, 21.48823577356381, 8.174118170771855, 33.74658417836842, 17.95641638953925),
    @(12, 11.15858320719366, 7.380286590078088, 11.00583599460683, 44.6471693777734, 66.34411907801265, 21.53850315236689, 8.171533344157647, 33.98374746715153, 18.015383087099),
    @(13, 11.16095458744504, 7.380246122028986, 11.00949277567637, 44.60222204979298, 66.25867445089727, 21.52763392113001, 8.172086066942493, 33.93278386463302, 18.00268684624352),
    @(14, 11.16877266505207, 7.380182757016339, 11.02146155726288, 44.45579895646335, 65.98000053362021, 21.49235176975694, 8.173903727669998, 33.76614672860651, 17.96126903062496),
    @(15, 11.17362627244343, 7.38019656758244, 11.02882603470634, 44.36623699274169, 65.80929628168326, 21.47086756175356, 8.175028718850182, 33.66374682948024, 17.93589058183993),
    @(16, 11.20242893963608, 7.381078198425432, 11.07154017862658, 43.85507748313351, 64.83123372738342, 21.34971956694471, 8.181654106963322, 33.07202555136449, 17.79035687007953),
    @(17, 11.22097547155113, 7.382330605024216, 11.09820192264313, 43.54351720130379, 64.23171586640323, 21.27717839652876, 8.185878288473159, 32.70480876049022, 17.70102617902899),
    @(18, 11.23196323461273, 7.383310626769458, 11.11370606088363, 43.36512185459198, 63.88720539774423, 21.23611690089419, 8.188366790768132, 32.49211437406295, 17.64964249827717),
    @(19, 11.23573832222688, 7.383686868098563, 11.11898458866517, 43.30486472925415, 63.77062747442373, 21.22232869939103, 8.189219480838226, 32.41985202489531, 17.63224582901911),
    @(20, 11.21896796769746, 7.382170367672735, 11.09534625840392, 43.57660135928571, 64.29550591955892, 21.2848321708701, 8.185422528638579, 32.74405446575027, 17.71053624899562),
    @(21, 11.16646278858275, 7.380190392486736, 11.0179391099356, 44.49878097392759, 66.06185629436338, 21.5026885499364, 8.173367415679314, 33.81516117720754, 17.97343637382503),
    @(22, 11.13504287230498, 7.381246809657621, 10.96883926730091, 45.10730770008772, 67.21625239987461, 21.65078342203869, 8.166009306518974, 34.50060711889392, 18.14490269235032),
    @(23, 11.15154524660717, 7.380466121469123, 10.99490914293767, 44.78205678331028, 66.60026403031542, 21.57122900429423, 8.169889006686059, 34.13616984170072, 18.05343590806848),
    @(24, 11.21987454800521, 7.382242001578162, 11.09663675600897, 43.56164173880772, 64.26666590513173, 21.28136989348891, 8.185628390845487, 32.72631637701831, 17.70623682454879),
    @(25, 11.30579959966355, 7.394060050458531, 11.21283690495735, 42.27684458667463, 61.76344240946423, 20.99411338191643, 8.204858489139632, 31.1499627952207, 17.33181221798671)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value  = $row[1]   # C
    $ws.Cells.Item($r, 4).Value  = $row[2]   # D
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G
    $ws.Cells.Item($r, 8).Value  = $row[6]   # H
    $ws.Cells.Item($r, 12).Value = $row[7]   # L
    $ws.Cells.Item($r, 13).Value = $row[8]   # M
    $ws.Cells.Item($r, 14).Value = $row[9]   # N
}

Write-Output "Updated loading_percent values for rows 2-25 (case with 380 kV)"
